$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 1821425.861256134
$ws.Range("B10").Value = 1860212.117059181
$ws.Range("B11").Value = 1895688.046766983
$ws.Range("B12").Value = 1916925.20873892
$ws.Range("B13").Value = 1917946.270817918
$ws.Range("B14").Value = 1908660.572850714
$ws.Range("B15").Value = 1885280.134636559
$ws.Range("B16").Value = 1845537.119791201
$ws.Range("B17").Value = 1794449.081568098
$ws.Range("B18").Value = 1739089.021255066
$ws.Range("B19").Value = 1676151.514983577
$ws.Range("B20").Value = 1592970.386680482
$ws.Range("B21").Value = 1501051.380243554
$ws.Range("B22").Value = 1407854.059279555
$ws.Range("B23").Value = 1304176.832156197
$ws.Range("B24").Value = 1194085.255215036
